$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: name and card number
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 02.05.2025"

# Row 6
$ws.Range("B6").Value = "04.05."
$ws.Range("C6").Value = "05.05."
$ws.Range("D6").Value = "KARTENZ./04.05 EDEKA RO"
$ws.Range("E6").Value = "146,75-"

# Row 7
$ws.Range("B7").Value = "07.05."
$ws.Range("C7").Value = "08.05."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 49394708"
$ws.Range("E7").Value = "86,90-"

# Row 8
$ws.Range("B8").Value = "10.05."
$ws.Range("C8").Value = "11.05."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,58-"

# Row 12: closing balance date / amount
$ws.Range("D12").Value = "KONTOSTAND AM 15.05.2025"
$ws.Range("E12").Value = "258,23-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 25.05.2025"
